$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 contains the candidate record; update the changed fields.
$ws.Range("A2").Value = "tCkfV986"     # Client Id
$ws.Range("B2").Value = 23091346       # Candidate ID
$ws.Range("C2").Value = "sjfmbpc53"    # User Name
$ws.Range("D2").Value = "Jk#5T%8h"     # Exam Password
$ws.Range("F2").Value = "uCWxYulc"     # First Name
$ws.Range("G2").Value = "qsZd"         # Last Name
